$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 54.66666666666666
$ws.Range("G4").Value = 58.66666666666666

$ws.Range("F5").Value = 30.66666666666666
$ws.Range("G5").Value = 33.33333333333333

$ws.Range("F6").Value = 16
$ws.Range("G6").Value = 12.66666666666667
